$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 137 and 138 (columns B:AC), keep column A fixed ---
$r137 = $ws.Range("B137:AC137")
$r138 = $ws.Range("B138:AC138")
$v137 = $r137.Value()
$v138 = $r138.Value()
$r137.Value = $v138
$r138.Value = $v137

# --- Swap rows 143 and 145 (columns B:AC), keep column A fixed ---
$r143 = $ws.Range("B143:AC143")
$r145 = $ws.Range("B145:AC145")
$v143 = $r143.Value()
$v145 = $r145.Value()
$r143.Value = $v145
$r145.Value = $v143

# --- Append new rows 197-202 ---

# Row 197
$ws.Range("A2").Copy()
$ws.Range("A197").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E197").PasteSpecial(-4122)
$ws.Cells.Item(197,1).Value = 195
$ws.Cells.Item(197,2).Value = 7609197
$ws.Cells.Item(197,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(197,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(197,5).Value = 45362.8125
$ws.Cells.Item(197,6).Value = "Libertad Asuncion"
$ws.Cells.Item(197,7).Value = "Sportivo Ameliano"
$ws.Cells.Item(197,8).Value = 4
$ws.Cells.Item(197,9).Value = 1
$ws.Cells.Item(197,10).Value = "H"
$ws.Cells.Item(197,11).Value = 1.5
$ws.Cells.Item(197,12).Value = 4
$ws.Cells.Item(197,13).Value = 6.5
$ws.Cells.Item(197,14).Value = 1.45
$ws.Cells.Item(197,15).Value = 4
$ws.Cells.Item(197,16).Value = 6.5
$ws.Cells.Item(197,17).Value = -1
$ws.Cells.Item(197,18).Value = 1.8
$ws.Cells.Item(197,19).Value = 2
$ws.Cells.Item(197,20).Value = 2.75
$ws.Cells.Item(197,21).Value = 1.975
$ws.Cells.Item(197,22).Value = 1.825
$ws.Cells.Item(197,23).Value = 0.45
$ws.Cells.Item(197,24).Value = -1
$ws.Cells.Item(197,25).Value = -1
$ws.Cells.Item(197,26).Value = 0.8
$ws.Cells.Item(197,27).Value = -1
$ws.Cells.Item(197,28).Value = 0.9750000000000001
$ws.Cells.Item(197,29).Value = -1

# Row 198
$ws.Range("A2").Copy()
$ws.Range("A198").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E198").PasteSpecial(-4122)
$ws.Cells.Item(198,1).Value = 196
$ws.Cells.Item(198,2).Value = 7609198
$ws.Cells.Item(198,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(198,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(198,5).Value = 45366.8125
$ws.Cells.Item(198,6).Value = "Sportivo Luqueno"
$ws.Cells.Item(198,7).Value = "Cerro Porteno"
$ws.Cells.Item(198,8).Value = 0
$ws.Cells.Item(198,9).Value = 1
$ws.Cells.Item(198,10).Value = "A"
$ws.Cells.Item(198,11).Value = 3.8
$ws.Cells.Item(198,12).Value = 3.3
$ws.Cells.Item(198,13).Value = 1.85
$ws.Cells.Item(198,14).Value = 3.6
$ws.Cells.Item(198,15).Value = 3.5
$ws.Cells.Item(198,16).Value = 1.909
$ws.Cells.Item(198,17).Value = 0.5
$ws.Cells.Item(198,18).Value = 1.85
$ws.Cells.Item(198,19).Value = 1.95
$ws.Cells.Item(198,20).Value = 2.5
$ws.Cells.Item(198,21).Value = 1.925
$ws.Cells.Item(198,22).Value = 1.875
$ws.Cells.Item(198,23).Value = -1
$ws.Cells.Item(198,24).Value = -1
$ws.Cells.Item(198,25).Value = 0.909
$ws.Cells.Item(198,26).Value = -1
$ws.Cells.Item(198,27).Value = 0.95
$ws.Cells.Item(198,28).Value = -1
$ws.Cells.Item(198,29).Value = 0.875

# Row 199
$ws.Range("A2").Copy()
$ws.Range("A199").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E199").PasteSpecial(-4122)
$ws.Cells.Item(199,1).Value = 197
$ws.Cells.Item(199,2).Value = 7959257
$ws.Cells.Item(199,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(199,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(199,5).Value = 45367.77777777778
$ws.Cells.Item(199,6).Value = "Olimpia Asuncion"
$ws.Cells.Item(199,7).Value = "2 de Mayo"
$ws.Cells.Item(199,8).Value = 2
$ws.Cells.Item(199,9).Value = 1
$ws.Cells.Item(199,10).Value = "H"
$ws.Cells.Item(199,11).Value = 1.666
$ws.Cells.Item(199,12).Value = 3.75
$ws.Cells.Item(199,13).Value = 4.333
$ws.Cells.Item(199,14).Value = 1.909
$ws.Cells.Item(199,15).Value = 3.6
$ws.Cells.Item(199,16).Value = 3.5
$ws.Cells.Item(199,17).Value = -0.5
$ws.Cells.Item(199,18).Value = 1.95
$ws.Cells.Item(199,19).Value = 1.85
$ws.Cells.Item(199,20).Value = 2.25
$ws.Cells.Item(199,21).Value = 1.85
$ws.Cells.Item(199,22).Value = 1.95
$ws.Cells.Item(199,23).Value = 0.909
$ws.Cells.Item(199,24).Value = -1
$ws.Cells.Item(199,25).Value = -1
$ws.Cells.Item(199,26).Value = 0.95
$ws.Cells.Item(199,27).Value = -1
$ws.Cells.Item(199,28).Value = 0.8500000000000001
$ws.Cells.Item(199,29).Value = -1

# Row 200
$ws.Range("A2").Copy()
$ws.Range("A200").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E200").PasteSpecial(-4122)
$ws.Cells.Item(200,1).Value = 198
$ws.Cells.Item(200,2).Value = 7609200
$ws.Cells.Item(200,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(200,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(200,5).Value = 45368.77777777778
$ws.Cells.Item(200,6).Value = "Nacional Asuncion"
$ws.Cells.Item(200,7).Value = "Tacuary"
$ws.Cells.Item(200,11).Value = 1.909
$ws.Cells.Item(200,12).Value = 3.4
$ws.Cells.Item(200,13).Value = 3.6
$ws.Cells.Item(200,14).Value = 2
$ws.Cells.Item(200,15).Value = 3.3
$ws.Cells.Item(200,16).Value = 3.3
$ws.Cells.Item(200,17).Value = -0.25
$ws.Cells.Item(200,18).Value = 1.925
$ws.Cells.Item(200,19).Value = 1.875
$ws.Cells.Item(200,20).Value = 2.25
$ws.Cells.Item(200,21).Value = 1.95
$ws.Cells.Item(200,22).Value = 1.85
$ws.Cells.Item(200,23).Value = 0
$ws.Cells.Item(200,24).Value = 0
$ws.Cells.Item(200,25).Value = 0
$ws.Cells.Item(200,26).Value = 0
$ws.Cells.Item(200,27).Value = 0

# Row 201
$ws.Range("A2").Copy()
$ws.Range("A201").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E201").PasteSpecial(-4122)
$ws.Cells.Item(201,1).Value = 199
$ws.Cells.Item(201,2).Value = 7609143
$ws.Cells.Item(201,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(201,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(201,5).Value = 45368.875
$ws.Cells.Item(201,6).Value = "Guarani Asuncion"
$ws.Cells.Item(201,7).Value = "Sportivo Trinidense"
$ws.Cells.Item(201,11).Value = 1.909
$ws.Cells.Item(201,12).Value = 3.4
$ws.Cells.Item(201,13).Value = 3.5
$ws.Cells.Item(201,14).Value = 1.909
$ws.Cells.Item(201,15).Value = 3.4
$ws.Cells.Item(201,16).Value = 3.6
$ws.Cells.Item(201,17).Value = -0.5
$ws.Cells.Item(201,18).Value = 1.9
$ws.Cells.Item(201,19).Value = 1.9
$ws.Cells.Item(201,20).Value = 2.5
$ws.Cells.Item(201,21).Value = 2
$ws.Cells.Item(201,22).Value = 1.8
$ws.Cells.Item(201,23).Value = 0
$ws.Cells.Item(201,24).Value = 0
$ws.Cells.Item(201,25).Value = 0
$ws.Cells.Item(201,26).Value = 0
$ws.Cells.Item(201,27).Value = 0

# Row 202
$ws.Range("A2").Copy()
$ws.Range("A202").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E202").PasteSpecial(-4122)
$ws.Cells.Item(202,1).Value = 200
$ws.Cells.Item(202,2).Value = 7609199
$ws.Cells.Item(202,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(202,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(202,5).Value = 45369.8125
$ws.Cells.Item(202,6).Value = "Sportivo Ameliano"
$ws.Cells.Item(202,7).Value = "Sol de America"
$ws.Cells.Item(202,11).Value = 1.909
$ws.Cells.Item(202,12).Value = 3.4
$ws.Cells.Item(202,13).Value = 3.6
$ws.Cells.Item(202,14).Value = 2.1
$ws.Cells.Item(202,15).Value = 3.3
$ws.Cells.Item(202,16).Value = 3.1
$ws.Cells.Item(202,17).Value = -0.25
$ws.Cells.Item(202,18).Value = 1.875
$ws.Cells.Item(202,19).Value = 1.925
$ws.Cells.Item(202,20).Value = 2.25
$ws.Cells.Item(202,21).Value = 1.775
$ws.Cells.Item(202,22).Value = 2.025
$ws.Cells.Item(202,23).Value = 0
$ws.Cells.Item(202,24).Value = 0
$ws.Cells.Item(202,25).Value = 0
$ws.Cells.Item(202,26).Value = 0
$ws.Cells.Item(202,27).Value = 0
